# Fixed naive component forecaster bug - Presentation state 11.02.
# The matured-error series for each forecast origin (row) was missing its
# own-quarter (Q0) naive error: every Q0..Q8 value had been written one
# column too far right (into Q1..Q9), so Q0 was blank and the true Q9 spilled
# off the end. This rebuilds B2:K20 with Q0 restored in column B and the
# rest of each row shifted back into its correct Q-column, dropping the
# stale trailing value that no longer belongs to the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 19,10

# Row 2 (A2)
$data[0,0] = 0.6014450472570072
$data[0,1] = -1.19942097753403
$data[0,2] = -2.158355674041143
$data[0,3] = 0.6753625508398458
$data[0,4] = -0.9937421128663182
$data[0,5] = 0.02643032487692459
$data[0,6] = -1.687359132022387
$data[0,7] = 1.140996241576585
$data[0,8] = -1.094336358289534
$data[0,9] = 0.2569986996281282
# Row 3 (A3)
$data[1,0] = -1.308035509750171
$data[1,1] = -2.266970206257284
$data[1,2] = 0.5667480186237051
$data[1,3] = -1.102356645082459
$data[1,4] = -0.08218420733921622
$data[1,5] = -1.795973664238528
$data[1,6] = 1.032381709360444
$data[1,7] = -1.202950890505675
$data[1,8] = 0.1483841674119874
$data[1,9] = -0.4730537624049144
# Row 4 (A4)
$data[2,0] = -2.058868485289545
$data[2,1] = 0.774849739591444
$data[2,2] = -0.89425492411472
$data[2,3] = 0.1259175136285228
$data[2,4] = -1.587871943270789
$data[2,5] = 1.240483430328183
$data[2,6] = -0.994849169537936
$data[2,7] = 0.3564858883797264
$data[2,8] = -0.2649520414371754
$data[2,9] = 0.2356086557536556
# Row 5 (A5)
$data[3,0] = 0.8913479099652445
$data[3,1] = -0.7777567537409195
$data[3,2] = 0.2424156840023232
$data[3,3] = -1.471373772896988
$data[3,4] = 1.356981600701984
$data[3,5] = -0.8783509991641355
$data[3,6] = 0.4729840587535268
$data[3,7] = -0.148453871063375
$data[3,8] = 0.352106826127456
$data[3,9] = -0.3158012462186854
# Row 6 (A6)
$data[4,0] = -0.7322633397437844
$data[4,1] = 0.2879090979994584
$data[4,2] = -1.425880358899853
$data[4,3] = 1.402475014699119
$data[4,4] = -0.8328575851670005
$data[4,5] = 0.5184774727506619
$data[4,6] = -0.1029604570662399
$data[4,7] = 0.3976002401245912
$data[4,8] = -0.2703078322215502
$data[4,9] = 0.1586931430164528
# Row 7 (A7)
$data[5,0] = 0.2703549766394939
$data[5,1] = -1.443434480259818
$data[5,2] = 1.384920893339154
$data[5,3] = -0.8504117065269649
$data[5,4] = 0.5009233513906975
$data[5,5] = -0.1205145784262043
$data[5,6] = 0.3800461187646267
$data[5,7] = -0.2878619535815147
$data[5,8] = 0.1411390216564884
$data[5,9] = 0.1988105702346985
# Row 8 (A8)
$data[6,0] = -1.355327161308811
$data[6,1] = 1.473028212290161
$data[6,2] = -0.7623043875759586
$data[6,3] = 0.5890306703417038
$data[6,4] = -0.0324072594751981
$data[6,5] = 0.4681534377156329
$data[6,6] = -0.1997546346305085
$data[6,7] = 0.2292463406074946
$data[6,8] = 0.2869178891857047
$data[6,9] = 0.4108842600239903
# Row 9 (A9)
$data[7,0] = 1.651602845777944
$data[7,1] = -0.5837297540881751
$data[7,2] = 0.7676053038294873
$data[7,3] = 0.1461673740125855
$data[7,4] = 0.6467280712034165
$data[7,5] = -0.02118000114272489
$data[7,6] = 0.4078209740952782
$data[7,7] = 0.4654925226734883
$data[7,8] = 0.5894588935117738
$data[7,9] = -0.303959229723018
# Row 10 (A10)
$data[8,0] = 0.3282974736644749
$data[8,1] = 1.679632531582137
$data[8,2] = 1.058194601765235
$data[8,3] = 1.558755298956066
$data[8,4] = 0.8908472266099251
$data[8,5] = 1.319848201847928
$data[8,6] = 1.377519750426138
$data[8,7] = 1.501486121264424
$data[8,8] = 0.608067998029632
$data[8,9] = 1.378198724973394
# Row 11 (A11)
$data[9,0] = 0.7356582956163805
$data[9,1] = 0.1142203657994787
$data[9,2] = 0.6147810629903097
$data[9,3] = -0.0531270093558317
$data[9,4] = 0.3758739658821714
$data[9,5] = 0.4335455144603815
$data[9,6] = 0.557511885298667
$data[9,7] = -0.3359062379361248
$data[9,8] = 0.4342244890076376
$data[9,9] = 0.1683237681281231
# Row 12 (A12)
$data[10,0] = 0.1181882633125878
$data[10,1] = 0.6187489605034189
$data[10,2] = -0.04915911184272259
$data[10,3] = 0.3798418633952805
$data[10,4] = 0.4375134119734906
$data[10,5] = 0.5614797828117761
$data[10,6] = -0.3319383404230157
$data[10,7] = 0.4381923865207467
$data[10,8] = 0.1722916656412322
$data[10,9] = $null
# Row 13 (A13)
$data[11,0] = 0.7543890506736601
$data[11,1] = 0.08648097832751878
$data[11,2] = 0.5154819535655218
$data[11,3] = 0.573153502143732
$data[11,4] = 0.6971198729820175
$data[11,5] = -0.1962982502527744
$data[11,6] = 0.5738324766909881
$data[11,7] = 0.3079317558114735
$data[11,8] = $null
$data[11,9] = $null
# Row 14 (A14)
$data[12,0] = -0.1543252035281459
$data[12,1] = 0.2746757717098572
$data[12,2] = 0.3323473202880673
$data[12,3] = 0.4563136911263528
$data[12,4] = -0.4371044321084391
$data[12,5] = 0.3330262948353234
$data[12,6] = 0.06712557395580883
$data[12,7] = $null
$data[12,8] = $null
$data[12,9] = $null
# Row 15 (A15)
$data[13,0] = 0.2293445564577608
$data[13,1] = 0.2870161050359709
$data[13,2] = 0.4109824758742565
$data[13,3] = -0.4824356473605354
$data[13,4] = 0.287695079583227
$data[13,5] = 0.02179435870371246
$data[13,6] = $null
$data[13,7] = $null
$data[13,8] = $null
$data[13,9] = $null
# Row 16 (A16)
$data[14,0] = 0.2201546830999171
$data[14,1] = 0.3441210539382026
$data[14,2] = -0.5492970692965893
$data[14,3] = 0.2208336576471732
$data[14,4] = -0.04506706323234141
$data[14,5] = $null
$data[14,6] = $null
$data[14,7] = $null
$data[14,8] = $null
$data[14,9] = $null
# Row 17 (A17)
$data[15,0] = 0.314534851581486
$data[15,1] = -0.5788832716533059
$data[15,2] = 0.1912474552904566
$data[15,3] = -0.07465326558905801
$data[15,4] = $null
$data[15,5] = $null
$data[15,6] = $null
$data[15,7] = $null
$data[15,8] = $null
$data[15,9] = $null
# Row 18 (A18)
$data[16,0] = -0.5970339283829468
$data[16,1] = 0.1730967985608157
$data[16,2] = -0.0928039223186989
$data[16,3] = $null
$data[16,4] = $null
$data[16,5] = $null
$data[16,6] = $null
$data[16,7] = $null
$data[16,8] = $null
$data[16,9] = $null
# Row 19 (A19)
$data[17,0] = 0.1550649743121164
$data[17,1] = -0.1108357465673982
$data[17,2] = $null
$data[17,3] = $null
$data[17,4] = $null
$data[17,5] = $null
$data[17,6] = $null
$data[17,7] = $null
$data[17,8] = $null
$data[17,9] = $null
# Row 20 (A20)
$data[18,0] = -0.1624199859130616
$data[18,1] = $null
$data[18,2] = $null
$data[18,3] = $null
$data[18,4] = $null
$data[18,5] = $null
$data[18,6] = $null
$data[18,7] = $null
$data[18,8] = $null
$data[18,9] = $null

$ws.Range("B2:K20").Value = $data
